$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "7-nov" column header (CG1), matching the style of the
# preceding date header cell (CF1)
$ws.Range("CG1").Value = "7-nov"
$ws.Range("CG1").NumberFormat = $ws.Range("CF1").NumberFormat

# Fill in the new column's data values (CG2:CG11), copying the number
# format / alignment of the corresponding CF cell in each row
$values = @(7, 10, 7, 9, 9, 9, 10, 11, 8, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $src = $ws.Cells.Item($row, 84)
    $dst = $ws.Cells.Item($row, 85)
    $dst.Value = $values[$i]
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.NumberFormat = $src.NumberFormat
}

# Move the active selection to CG8, matching the edited file
$ws.Range("CG8").Select()
